# Update the Galp-Galr2 LR-pair worksheet with new TPM-derived values.
# A new "ECs" sending/target cluster row is inserted (row 2), the
# remaining target-cluster rows shift down by one, and all numeric
# metric columns are refreshed with the newly recalculated TPM values.
# After the edit the sheet has 6 data rows (rows 2-7) instead of 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2: Target cluster = ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Galp"
$ws.Range("C2").Value = "Galr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.5
$ws.Range("G2").Value = 0.1314505
$ws.Range("H2").Value = 0.262901
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.5954815
$ws.Range("N2").Value = 1.190963
$ws.Range("O2").Value = 0.1769072473669437
$ws.Range("P2").Value = 0.1359316837037387
$ws.Range("Q2").Value = 0.07827634091575
$ws.Range("R2").Value = 0.313105363663
$ws.Range("S2").Value = 0.1769072473669437
$ws.Range("T2").Value = 0.1359316837037387

# Row 3: Target cluster = FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Galp"
$ws.Range("C3").Value = "Galr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.5
$ws.Range("G3").Value = 0.1314505
$ws.Range("H3").Value = 0.262901
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.2698953333333333
$ws.Range("N3").Value = 0.809686
$ws.Range("O3").Value = 0.08018123232574602
$ws.Range("P3").Value = 0.09241427420612175
$ws.Range("Q3").Value = 0.03547787651433333
$ws.Range("R3").Value = 0.212867259086
$ws.Range("S3").Value = 0.08018123232574602
$ws.Range("T3").Value = 0.09241427420612175

# Row 4: Target cluster = Inflammatory-Mac
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Galp"
$ws.Range("C4").Value = "Galr2"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = 0.1314505
$ws.Range("H4").Value = 0.262901
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1259886666666667
$ws.Range("N4").Value = 0.377966
$ws.Range("O4").Value = 0.03742905232056985
$ws.Range("P4").Value = 0.04313950539417875
$ws.Range("Q4").Value = 0.01656127322766666
$ws.Range("R4").Value = 0.099367639366
$ws.Range("S4").Value = 0.03742905232056985
$ws.Range("T4").Value = 0.04313950539417875

# Row 5: Target cluster = MuSCs
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Galp"
$ws.Range("C5").Value = "Galr2"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.5
$ws.Range("G5").Value = 0.1314505
$ws.Range("H5").Value = 0.262901
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.741235
$ws.Range("N5").Value = 1.48247
$ws.Range("O5").Value = 0.2202080895914256
$ws.Range("P5").Value = 0.169203109702217
$ws.Range("Q5").Value = 0.0974357113675
$ws.Range("R5").Value = 0.38974284547
$ws.Range("S5").Value = 0.2202080895914256
$ws.Range("T5").Value = 0.169203109702217

# Row 6: Target cluster = Neutrophils
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Galp"
$ws.Range("C6").Value = "Galr2"
$ws.Range("D6").Value = "Neutrophils"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.1314505
$ws.Range("H6").Value = 0.262901
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.177331
$ws.Range("N6").Value = 3.531993
$ws.Range("O6").Value = 0.3497646634694297
$ws.Range("P6").Value = 0.4031273476336537
$ws.Range("Q6").Value = 0.1547607486155
$ws.Range("R6").Value = 0.9285644916929999
$ws.Range("S6").Value = 0.3497646634694297
$ws.Range("T6").Value = 0.4031273476336537

# Row 7: Target cluster = Resolving-Mac
$ws.Range("A7").Value = "ECs"
$ws.Range("B7").Value = "Galp"
$ws.Range("C7").Value = "Galr2"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.1314505
$ws.Range("H7").Value = 0.262901
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 1
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.4561346666666666
$ws.Range("N7").Value = 1.368404
$ws.Range("O7").Value = 0.135509714925885
$ws.Range("P7").Value = 0.15618407936009
$ws.Range("Q7").Value = 0.05995913000066666
$ws.Range("R7").Value = 0.359754780004
$ws.Range("S7").Value = 0.135509714925885
$ws.Range("T7").Value = 0.15618407936009
